# Progress report 5 update for Worklogs(NEn408).xlsx
# - Updates the repo-link hyperlink cell text (E6) to the new repo URL
# - Appends 5 new worklog rows (20-24) covering Nov 9 - Nov 15, 2025
# - Moves the sheet selection to B3 (scrolled back to top)
# - Sets the print scale to 88%

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the repo link text shown in E6 (hyperlink target itself is left as-is) ---
$ws.Range("E6").Value = "https://github.com/narkmn/F2025_4495_050_Nen408"

# --- Add new rows 20-24, cloning formatting from the last existing row (19) ---
$ws.Range("B19:E19").Copy()
$ws.Range("B20:E24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 20
$ws.Range("B20").Value = "11/9/2025"
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = "Wordpress migration"
$ws.Range("E20").Value = "Improved migration code, because of database url problem, I created 3 migration scripts."
$ws.Rows("20").RowHeight = 28.8

# Row 21
$ws.Range("B21").Value = "11/9/2025"
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = "Learning"
$ws.Range("E21").Value = "Researched openai chatbot, and custom plugin of wordpress"
$ws.Rows("21").RowHeight = 28.8

# Row 22
$ws.Range("B22").Value = "11/13/2025"
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = "Learning"
$ws.Range("E22").Value = "Tested editing custom plugin in www.dataofattraction.com, however something is broken, now I cannot access use customize tool"
$ws.Rows("22").RowHeight = 43.2

# Row 23
$ws.Range("B23").Value = "11/14/2025"
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = "AI chatbot"
$ws.Range("E23").Value = "Created AI chatbot plugin in local wordpress host. It is working fine, but it work like Chatgpt, not personalized and related to topic"
$ws.Rows("23").RowHeight = 43.2

# Row 24
$ws.Range("B24").Value = "11/15/2025"
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = "AI chatbot"
$ws.Range("E24").Value = "added the plugin into healthcare.ca and localhost. And get user info and saved chat history in database. Added last 10 chat history along with all class topic since there is no way which class student in"
$ws.Rows("24").RowHeight = 57.6

# --- Print scale ---
$ws.PageSetup.Zoom = 88

# --- Selection / view: back to top-left, select B3 ---
$ws.Range("B3").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# --- Best-effort window geometry (mirrors the saved workbook window size/position) ---
try {
    $excel.ActiveWindow.Left = -108
    $excel.ActiveWindow.Top = -108
    $excel.ActiveWindow.Width = 23256
    $excel.ActiveWindow.Height = 12456
} catch {}

Write-Host "Worklog updated: hyperlink text refreshed, 5 rows added, view + print scale updated."
